$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$old = "dnasr281@gmail.com, System"
$new = "System, dnasr281@gmail.com"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text
    if ($val -eq $old) {
        $cell.Value = $new
    }
}
